$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.753.18"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.062.33"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.23%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "517.33"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.78"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.84%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("E9").Value = "  +2.57%  "
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("E11").Value = "  +3.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.588.16"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.31%  "
$ws.Range("E13").Value = "  +3.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.92"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.47%  "
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.787.76"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.066.77"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.09"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.95%  "
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.11"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "330.85"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.71"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.169"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.51%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  -3.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.37"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.19"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.12%  "
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("E31").Value = "  +3.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.73"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.69"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.49"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("E35").Value = "  +3.03%  "
$ws.Range("E36").Value = "  +1.50%  "
$ws.Range("E37").Value = "  +1.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0673"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.105.04"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.38%  "
$ws.Range("E40").Value = "  +3.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.53"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.254.13"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0258"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +9.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "20.58"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +6.06%  "
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.933"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.731"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +8.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "258.19"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +12.51%  "
